# NIT-9012564319.xlsx — "Elimna EC anteriores y se agregan nuevos, se
# modifica base de datos"
#
# The "Periodo Mora" column (E16:E22) is reordered so the most recent
# period (2309) is listed first and the oldest (2303) last, and the
# "Valor Mora" amounts in F16/F22 swap along with the row they now sit
# in (the 24000 outlier moves from the last period row to the first).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Periodo Mora (column E), rows 16-22: reverse the period order ---
$periods = @("2309", "2308", "2307", "2306", "2305", "2304", "2303")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# --- Valor Mora (column F): the 24000 value moves from the last row
#     (22) to the first (16); the rest stay at 40000 ---
$ws.Range("F16").Value = 24000
$ws.Range("F17").Value = 40000
$ws.Range("F18").Value = 40000
$ws.Range("F19").Value = 40000
$ws.Range("F20").Value = 40000
$ws.Range("F21").Value = 40000
$ws.Range("F22").Value = 40000

# --- Column widths (bestFit recalculated by Excel on save) ---
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666

$wb.Save()
